$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 187.04546
$ws.Range("I9").Value = 161.64706
$ws.Range("J9").Value = 273.4
$ws.Range("K9").Value = 161.64706
$ws.Range("L9").Value = 273.4
$ws.Range("M9").Value = 7.35293999999999
$ws.Range("N9").Value = -611.4

# Row 18
$ws.Range("H18").Value = 1527
$ws.Range("I18").Value = 1527
$ws.Range("K18").Value = 1527
$ws.Range("M18").Value = -1243

# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("N27").Value = 0

# Row 41
$ws.Range("H41").Value = 405.375
$ws.Range("I41").Value = 423.5
$ws.Range("K41").Value = 423.5
$ws.Range("M41").Value = 16.5

# Row 53
$ws.Range("H53").Value = 2369.4348
$ws.Range("I53").Value = 4124.5386
$ws.Range("K53").Value = 4124.5386
$ws.Range("M53").Value = -3487.5386

# Row 82
$ws.Range("H82").Value = 10001193
$ws.Range("I82").Value = 10001193
$ws.Range("K82").Value = 30003579
$ws.Range("M82").Value = -30003173

# Row 85
$ws.Range("H85").Value = 10001193
$ws.Range("I85").Value = 10001193
$ws.Range("K85").Value = 30003579
$ws.Range("M85").Value = -30002175

# Row 87
$ws.Range("H87").Value = 58894.445
$ws.Range("J87").Value = 58894.445
$ws.Range("L87").Value = 58894.445
$ws.Range("N87").Value = -61390.445

# Row 90
$ws.Range("H90").Value = 58894.445
$ws.Range("J90").Value = 58894.445
$ws.Range("L90").Value = 176683.335
$ws.Range("N90").Value = -189163.335

# Row 97
$ws.Range("H97").Value = 5223.3335
$ws.Range("J97").Value = 5223.3335
$ws.Range("L97").Value = 15670.0005
$ws.Range("N97").Value = -16662.0005

# Row 101
$ws.Range("H101").Value = 369
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()

# Row 118
$ws.Range("H118").Value = 4842.857
$ws.Range("I118").Value = 4842.857
$ws.Range("K118").Value = 14528.571
$ws.Range("M118").Value = -12871.571

# Row 132
$ws.Range("H132").Value = 4573.7427
$ws.Range("I132").Value = 2551.1614
$ws.Range("K132").Value = 7653.4842
$ws.Range("M132").Value = -5123.4842

# Row 137
$ws.Range("H137").Value = 2337.05
$ws.Range("I137").Value = 3274.75
$ws.Range("J137").Value = 1711.9166
$ws.Range("K137").Value = 9824.25
$ws.Range("L137").Value = 5135.7498
$ws.Range("M137").Value = -7274.25
$ws.Range("N137").Value = -10235.7498

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 3891.2307
$ws.Range("I5").Value = 1127.125
$ws.Range("J5").Value = 8313.799999999999
$ws.Range("K5").Value = 1127.125
$ws.Range("L5").Value = 8313.799999999999
$ws.Range("M5").Value = -1015.125
$ws.Range("N5").Value = -8537.799999999999

# Row 63
$ws.Range("H63").Value = 1475.2
$ws.Range("I63").Value = 1475.2
$ws.Range("K63").Value = 1475.2
$ws.Range("M63").Value = -789.2

# Row 66
$ws.Range("H66").Value = 1475.2
$ws.Range("I66").Value = 1475.2
$ws.Range("K66").Value = 7376
$ws.Range("M66").Value = -3944

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 3891.2307
$ws.Range("I4").Value = 1127.125
$ws.Range("J4").Value = 8313.799999999999
$ws.Range("K4").Value = 1127.125
$ws.Range("L4").Value = 8313.799999999999
$ws.Range("M4").Value = -1012.125
$ws.Range("N4").Value = -8543.799999999999

# Row 132
$ws.Range("H132").Value = 99995
$ws.Range("J132").Value = 99995
$ws.Range("L132").Value = 99995
$ws.Range("N132").Value = -110115

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4245.387
$ws.Range("I31").Value = 1075.4166
$ws.Range("J31").Value = 6247.4736
$ws.Range("K31").Value = 1075.4166
$ws.Range("L31").Value = 6247.4736
$ws.Range("M31").Value = -780.4166
$ws.Range("N31").Value = -6837.4736

# Row 34
$ws.Range("H34").Value = 4245.387
$ws.Range("I34").Value = 1075.4166
$ws.Range("J34").Value = 6247.4736
$ws.Range("K34").Value = 1075.4166
$ws.Range("L34").Value = 6247.4736
$ws.Range("M34").Value = -873.4166
$ws.Range("N34").Value = -6651.4736

# Row 41
$ws.Range("H41").Value = 44963.875
$ws.Range("J41").Value = 73699.5
$ws.Range("L41").Value = 73699.5
$ws.Range("N41").Value = -74555.5

# Row 134
$ws.Range("H134").Value = 2697.9656
$ws.Range("I134").Value = 2578.5715
$ws.Range("K134").Value = 7735.7145
$ws.Range("M134").Value = -5200.7145

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 444488.94
$ws.Range("I2").Value = 625022.5600000001
$ws.Range("J2").Value = 98.46154
$ws.Range("K2").Value = 3750135.36
$ws.Range("L2").Value = 590.76924
$ws.Range("M2").Value = -3750022.36
$ws.Range("N2").Value = -816.76924

# Row 34
$ws.Range("H34").Value = 1511.6364
$ws.Range("J34").Value = 2012.375
$ws.Range("L34").Value = 6037.125
$ws.Range("N34").Value = -6205.125

# Row 39
$ws.Range("H39").Value = 5620
$ws.Range("J39").Value = 5999.75
$ws.Range("L39").Value = 17999.25
$ws.Range("N39").Value = -18587.25

# Row 48
$ws.Range("H48").Value = 4059.2
$ws.Range("J48").Value = 4324
$ws.Range("L48").Value = 12972
$ws.Range("N48").Value = -13472

# Row 75
$ws.Range("H75").Value = 10015
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 10015
$ws.Range("K75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("M75").Value = 30045
$ws.Range("N75").Value = -32041

# Row 78
$ws.Range("H78").Value = 10015
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 10015
$ws.Range("K78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("M78").Value = 90135
$ws.Range("N78").Value = -100119

# Row 122
$ws.Range("H122").Value = 1129.1875
$ws.Range("J122").Value = 1131.9
$ws.Range("L122").Value = 10187.1
$ws.Range("N122").Value = -15087.1

# Row 126
$ws.Range("H126").Value = 26999
$ws.Range("J126").Value = 26999
$ws.Range("L126").Value = 80997
$ws.Range("N126").Value = -90877

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 9397.799999999999
$ws.Range("I132").Value = 9572.710999999999
$ws.Range("J132").Value = 9004.25
$ws.Range("K132").Value = 28718.133
$ws.Range("L132").Value = 27012.75
$ws.Range("M132").Value = -26188.133
$ws.Range("N132").Value = -32072.75

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 19235034
$ws.Range("I7").Value = 4503.579
$ws.Range("K7").Value = 4503.579
$ws.Range("M7").Value = -4391.579

# Row 61
$ws.Range("H61").Value = 3198.85
$ws.Range("I61").Value = 2840.8948
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 2840.8948
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -2638.8948
$ws.Range("N61").Value = -10404

# Row 113
$ws.Range("H113").Value = 3198.85
$ws.Range("I113").Value = 2840.8948
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 2840.8948
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -670.8948
$ws.Range("N113").Value = -14340

# Row 126
$ws.Range("H126").Value = 19235034
$ws.Range("I126").Value = 4503.579
$ws.Range("K126").Value = 13510.737
$ws.Range("M126").Value = -11040.737

# Row 132
$ws.Range("H132").Value = 6568.8203
$ws.Range("I132").Value = 6035.5483
$ws.Range("K132").Value = 18106.6449
$ws.Range("M132").Value = -15576.6449

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2031.5714
$ws.Range("I81").Value = 1340
$ws.Range("K81").Value = 2680
$ws.Range("M81").Value = -1619

# Row 84
$ws.Range("H84").Value = 2031.5714
$ws.Range("I84").Value = 1340
$ws.Range("K84").Value = 13400
$ws.Range("M84").Value = -8096
